$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date number format used by column D (style s="2")
$dateFormat = $ws.Range("D58").NumberFormat()

# --- 1) Update column D (Fecha) for rows 58-113 ---
# Each pair of rows (Primera/Segunda) is shifted to the date of the
# following pair; the first pair gets a brand-new date, and the date
# that "falls off" the end is appended as new rows 114-115 below.
$dUpdates = @(
    @(58,44512),
    @(59,44512),
    @(60,44160),
    @(61,44160),
    @(62,44274),
    @(63,44274),
    @(64,44222),
    @(65,44222),
    @(66,44334),
    @(67,44334),
    @(68,44405),
    @(69,44405),
    @(70,44350),
    @(71,44350),
    @(72,44278),
    @(73,44278),
    @(74,44272),
    @(75,44272),
    @(76,44341),
    @(77,44341),
    @(78,44308),
    @(79,44308),
    @(80,44398),
    @(81,44398),
    @(82,44490),
    @(83,44490),
    @(84,44316),
    @(85,44316),
    @(86,44453),
    @(87,44453),
    @(88,44280),
    @(89,44280),
    @(90,44442),
    @(91,44442),
    @(92,44476),
    @(93,44476),
    @(94,44166),
    @(95,44166),
    @(96,44292),
    @(97,44292),
    @(98,44306),
    @(99,44306),
    @(100,44469),
    @(101,44469),
    @(102,44425),
    @(103,44425),
    @(104,44168),
    @(105,44168),
    @(106,44400),
    @(107,44400),
    @(108,44330),
    @(109,44330),
    @(110,44217),
    @(111,44217),
    @(112,44383),
    @(113,44383)
)

foreach ($u in $dUpdates) {
    $ws.Cells.Item($u[0], 4).Value = $u[1]
}

# --- 2) Update column J (Volumen) where the Primera/Segunda pairs swapped ---
$jUpdates = @(
    @(70,200),
    @(71,100),
    @(72,300),
    @(73,150),
    @(88,200),
    @(89,100),
    @(90,300),
    @(91,150)
)

foreach ($u in $jUpdates) {
    $ws.Cells.Item($u[0], 10).Value = $u[1]
}

# --- 3) Append new rows 114 and 115 (the date that fell off the shift) ---
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 44433
$ws.Cells.Item(114, 4).NumberFormat = $dateFormat
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100112044
$ws.Cells.Item(114, 7).Value = "Perejil"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 200
$ws.Cells.Item(114, 11).Value = 600
$ws.Cells.Item(114, 12).Value = 700
$ws.Cells.Item(114, 13).Value = 650
$ws.Cells.Item(114, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(114, 15).Value = "Región de Ñuble"
$ws.Cells.Item(114, 16).Value = 650
$ws.Cells.Item(114, 17).Value = 1
$ws.Cells.Item(114, 18).Value = "Hortaliza"

$ws.Cells.Item(115, 1).Value = 11
$ws.Cells.Item(115, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(115, 3).Value = "Bíobío"
$ws.Cells.Item(115, 4).Value = 44433
$ws.Cells.Item(115, 4).NumberFormat = $dateFormat
$ws.Cells.Item(115, 5).Value = 8
$ws.Cells.Item(115, 6).Value = 100112044
$ws.Cells.Item(115, 7).Value = "Perejil"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Segunda"
$ws.Cells.Item(115, 10).Value = 100
$ws.Cells.Item(115, 11).Value = 500
$ws.Cells.Item(115, 12).Value = 500
$ws.Cells.Item(115, 13).Value = 500
$ws.Cells.Item(115, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(115, 15).Value = "Región de Ñuble"
$ws.Cells.Item(115, 16).Value = 500
$ws.Cells.Item(115, 17).Value = 1
$ws.Cells.Item(115, 18).Value = "Hortaliza"

Write-Output "Done"
Write-Output $ws.UsedRange.Address()
